$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "CarModel"
$ws.Range("C1").Value = "LicensePlate"
$ws.Range("D1").Value = "ManufacturingYear"
$ws.Range("E1").Value = "Features"

# --- Id column ---
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# --- CarModel / LicensePlate columns (row by row) ---
$ws.Range("B2").Value = "CarModelA"
$ws.Range("C2").Value = "LicensePlate1"
$ws.Range("B3").Value = "CarModelB"
$ws.Range("C3").Value = "LicensePlate2"
$ws.Range("B4").Value = "CarModelC"
$ws.Range("C4").Value = "LicensePlate3"

# --- ManufacturingYear column ---
$ws.Range("D2").Value = 1985
$ws.Range("D3").Value = 1995
$ws.Range("D4").Value = 1992
$ws.Range("D4").Style = "Normal"

# --- Features column ---
$ws.Range("E2").Value = "AirConditioning`nPowerSteering"
$ws.Range("E2").WrapText = $true
$ws.Range("E3").Value = "PowerSteering`nBucketSeats"
$ws.Range("E3").WrapText = $true
$ws.Range("E4").Value = "AirConditioning`nBucketSeats"
$ws.Range("E4").WrapText = $true

# --- Column widths (approximate best-fit display widths) ---
$ws.Columns.Item(1).ColumnWidth = 1.833333333333333
$ws.Columns.Item(2).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 12.0
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 13.833333333333332

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("J9").Select() | Out-Null
